$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the greeting text for rule R10 (row 8, column E) from "Good Morning"
# to "GIT UPDATE".
$ws.Range("E8").Value = "GIT UPDATE"

# Make E8 the active/selected cell on the sheet.
$ws.Activate()
$ws.Range("E8").Select()
